$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: description ---
$g2 = @"
<get>
  <filter>
    <components xmlns="http://openconfig.net/yang/platform">
      <component>
        <name>Waveserver-Ai</name>
        <state>
          <description></description>
        </state>
      </component>
    </components>
  </filter>
</get>
"@
$ws.Range("G2").Value = $g2

$j2 = @"
<?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:1924bced-b91c-472c-8c3c-30d51b7bef9e"
 xmlns:ncx="http://netconfcentral.org/ns/yuma-ncx"
 ncx:last-modified="2020-10-07T13:51:28Z" ncx:etag="814"
 xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
 <data>
  <components xmlns="http://openconfig.net/yang/platform">
   <component>
    <name>Waveserver-Ai</name>
    <state>
     <description>Waveserver Ai Chassis 3-slot, 1RU</description>
    </state>
   </component>
  </components>
 </data>
</rpc-reply>
"@
$ws.Range("J2").Value = $j2

# --- Row 3: hardware-version ---
$g3 = @"
<get>
  <filter>
    <components xmlns="http://openconfig.net/yang/platform">
      <component>
        <name>Waveserver-Ai</name>
        <state>
          <hardware-version></hardware-version>
        </state>
      </component>
    </components>
  </filter>
</get>
"@
$ws.Range("G3").Value = $g3

$j3 = @"
<?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:f6bb4233-a6c3-495a-8e64-331b2883593a"
 xmlns:ncx="http://netconfcentral.org/ns/yuma-ncx"
 ncx:last-modified="2020-10-07T13:51:28Z" ncx:etag="814"
 xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
 <data>
  <components xmlns="http://openconfig.net/yang/platform">
   <component>
    <name>Waveserver-Ai</name>
    <state>
     <hardware-version>001</hardware-version>
    </state>
   </component>
  </components>
 </data>
</rpc-reply>
"@
$ws.Range("J3").Value = $j3

# --- Row 4: id ---
$g4 = @"
<get>
  <filter>
    <components xmlns="http://openconfig.net/yang/platform">
      <component>
        <name>Waveserver-Ai</name>
        <state>
          <id></id>
        </state>
      </component>
    </components>
  </filter>
</get>
"@
$ws.Range("G4").Value = $g4

$j4 = @"
<?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:aa1f2b73-992f-4608-8b3a-2c44ceafe004"
 xmlns:ncx="http://netconfcentral.org/ns/yuma-ncx"
 ncx:last-modified="2020-10-07T13:51:28Z" ncx:etag="814"
 xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
 <data>
  <components xmlns="http://openconfig.net/yang/platform">
   <component>
    <name>Waveserver-Ai</name>
    <state>
     <id>Waveserver Ai Chassis</id>
    </state>
   </component>
  </components>
 </data>
</rpc-reply>
"@
$ws.Range("J4").Value = $j4

# --- Row 5: location ---
$g5 = @"
<get>
  <filter>
    <components xmlns="http://openconfig.net/yang/platform">
      <component>
        <name>Waveserver-Ai</name>
        <state>
          <location></location>
        </state>
      </component>
    </components>
  </filter>
</get>
"@
$ws.Range("G5").Value = $g5

$j5 = @"
<?xml version="1.0" encoding="UTF-8"?>
<rpc-reply message-id="urn:uuid:393b75be-0861-4ec9-b887-37269807287d"
 xmlns:ncx="http://netconfcentral.org/ns/yuma-ncx"
 ncx:last-modified="2020-10-07T13:51:28Z" ncx:etag="814"
 xmlns="urn:ietf:params:xml:ns:netconf:base:1.0">
 <data>
  <components xmlns="http://openconfig.net/yang/platform">
   <component>
    <name>Waveserver-Ai</name>
    <state>
    </state>
   </component>
  </components>
 </data>
</rpc-reply>
"@
$ws.Range("J5").Value = $j5

$wb.Save()
